# Update the catalog: rename header labels and fill in the per-study id
# values that were missing on the "follow-up" sub-rows of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: shorten "tx" -> "t" and "baseline mean" -> "y"
$ws.Range("D1").Value = "t"
$ws.Range("E1").Value = "y"

# Fill in the study "id" on the second (follow-up) row of each study pair,
# mirroring the id already present on the first row of the pair.
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C5").Value = $ws.Range("C4").Value2
$ws.Range("C7").Value = $ws.Range("C6").Value2
$ws.Range("C9").Value = $ws.Range("C8").Value2
$ws.Range("C11").Value = $ws.Range("C10").Value2
$ws.Range("C13").Value = $ws.Range("C12").Value2
$ws.Range("C15").Value = $ws.Range("C14").Value2
$ws.Range("C17").Value = $ws.Range("C16").Value2
$ws.Range("C19").Value = $ws.Range("C18").Value2
$ws.Range("C21").Value = $ws.Range("C20").Value2
$ws.Range("C23").Value = $ws.Range("C22").Value2
$ws.Range("C25").Value = $ws.Range("C24").Value2

# Move the active selection, as left by the editor
$ws.Range("L20").Select() | Out-Null
